$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1001
$ws.Range("I5").Value = 1001
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1001
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -886
$ws.Range("H87").Value = 15602.921
$ws.Range("J87").Value = 15602.921
$ws.Range("L87").Value = 15602.921
$ws.Range("N87").Value = -18098.921
$ws.Range("H90").Value = 15602.921
$ws.Range("J90").Value = 15602.921
$ws.Range("L90").Value = 46808.763
$ws.Range("N90").Value = -59288.763
$ws.Range("H100").Value = 3192.9167
$ws.Range("I100").Value = 2696.6667
$ws.Range("J100").Value = 6666.6665
$ws.Range("K100").Value = 2696.6667
$ws.Range("L100").Value = 6666.6665
$ws.Range("M100").Value = -2155.6667
$ws.Range("N100").Value = -7748.6665
$ws.Range("H137").Value = 1012.7568
$ws.Range("I137").Value = 1014.4231
$ws.Range("J137").Value = 1008.8182
$ws.Range("K137").Value = 3043.2693
$ws.Range("L137").Value = 3026.4546
$ws.Range("M137").Value = -493.2692999999999
$ws.Range("N137").Value = -8126.4546
$ws.Range("H138").Value = 3663.55
$ws.Range("I138").Value = 2806.7
$ws.Range("J138").Value = 4030.7715
$ws.Range("K138").Value = 8420.099999999999
$ws.Range("L138").Value = 12092.3145
$ws.Range("M138").Value = -3280.099999999999
$ws.Range("N138").Value = -22372.3145
$ws.Range("H141").Value = 2134.3428
$ws.Range("I141").Value = 1051.5758
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 3154.7274
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = 2025.2726
$ws.Range("N141").Value = -70360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4598.12
$ws.Range("I32").Value = 4598.12
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4598.12
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4311.12
$ws.Range("H88").Value = 1604.2727
$ws.Range("I88").Value = 1200
$ws.Range("J88").Value = 1755.875
$ws.Range("K88").Value = 1200
$ws.Range("L88").Value = 1755.875
$ws.Range("M88").Value = -794
$ws.Range("N88").Value = -2567.875
$ws.Range("H91").Value = 1604.2727
$ws.Range("I91").Value = 1200
$ws.Range("J91").Value = 1755.875
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 1755.875
$ws.Range("M91").Value = 204
$ws.Range("N91").Value = -4563.875
$ws.Range("H133").Value = 63916.812
$ws.Range("J133").Value = 63916.812
$ws.Range("L133").Value = 63916.812
$ws.Range("N133").Value = -68976.81200000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 195.5
$ws.Range("I22").Value = 85.181816
$ws.Range("K22").Value = 85.181816
$ws.Range("M22").Value = 87.818184

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4116550
$ws.Range("I22").Value = 5292593
$ws.Range("K22").Value = 5292593
$ws.Range("M22").Value = -5292243
$ws.Range("H58").Value = 2541.7026
$ws.Range("I58").Value = 1527.1578
$ws.Range("J58").Value = 3612.611
$ws.Range("K58").Value = 1527.1578
$ws.Range("L58").Value = 3612.611
$ws.Range("M58").Value = -1324.1578
$ws.Range("N58").Value = -4018.611
$ws.Range("H115").Value = 32000
$ws.Range("J115").Value = 32000
$ws.Range("L115").Value = 32000
$ws.Range("N115").Value = -34350
$ws.Range("H136").Value = 2541.7026
$ws.Range("I136").Value = 1527.1578
$ws.Range("J136").Value = 3612.611
$ws.Range("K136").Value = 4581.4734
$ws.Range("L136").Value = 10837.833
$ws.Range("M136").Value = -2031.4734
$ws.Range("N136").Value = -15937.833

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 425897.5
$ws.Range("I107").Value = 473.25
$ws.Range("J107").Value = 567705.5600000001
$ws.Range("K107").Value = 1419.75
$ws.Range("L107").Value = 1703116.68
$ws.Range("M107").Value = 500.25
$ws.Range("N107").Value = -1706956.68

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1979.6222
$ws.Range("I132").Value = 1554.5862
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 4663.7586
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -2133.7586
$ws.Range("N132").Value = -13310

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 59083.945
$ws.Range("I7").Value = 94009.63
$ws.Range("J7").Value = 4200.7144
$ws.Range("K7").Value = 94009.63
$ws.Range("L7").Value = 4200.7144
$ws.Range("M7").Value = -93897.63
$ws.Range("N7").Value = -4424.7144
$ws.Range("H22").Value = 1013.8
$ws.Range("I22").Value = 1112.75
$ws.Range("J22").Value = 947.8333
$ws.Range("K22").Value = 1112.75
$ws.Range("L22").Value = 947.8333
$ws.Range("M22").Value = -817.75
$ws.Range("N22").Value = -1537.8333
$ws.Range("H27").Value = 1013.8
$ws.Range("I27").Value = 1112.75
$ws.Range("J27").Value = 947.8333
$ws.Range("K27").Value = 1112.75
$ws.Range("L27").Value = 947.8333
$ws.Range("M27").Value = -1005.75
$ws.Range("N27").Value = -1161.8333
$ws.Range("H46").Value = 2320.2354
$ws.Range("I46").Value = 1844.4
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1844.4
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1656.4
$ws.Range("N46").Value = -3376
$ws.Range("H68").Value = 1654.6428
$ws.Range("I68").Value = 1146
$ws.Range("J68").Value = 1739.4166
$ws.Range("K68").Value = 1146
$ws.Range("L68").Value = 1739.4166
$ws.Range("M68").Value = -397
$ws.Range("N68").Value = -3237.4166
$ws.Range("H71").Value = 1654.6428
$ws.Range("I71").Value = 1146
$ws.Range("J71").Value = 1739.4166
$ws.Range("K71").Value = 5730
$ws.Range("L71").Value = 8697.083000000001
$ws.Range("M71").Value = -1986
$ws.Range("N71").Value = -16185.083
$ws.Range("H126").Value = 59083.945
$ws.Range("I126").Value = 94009.63
$ws.Range("J126").Value = 4200.7144
$ws.Range("K126").Value = 282028.89
$ws.Range("L126").Value = 12602.1432
$ws.Range("M126").Value = -279558.89
$ws.Range("N126").Value = -17542.1432
$ws.Range("H127").Value = 49460.555
$ws.Range("J127").Value = 49460.555
$ws.Range("L127").Value = 49460.555
$ws.Range("N127").Value = -59380.555
$ws.Range("H130").Value = 41762.332
$ws.Range("J130").Value = 41762.332
$ws.Range("L130").Value = 41762.332
$ws.Range("N130").Value = -51802.332
$ws.Range("H136").Value = 12823526
$ws.Range("I136").Value = 3096.7058
$ws.Range("J136").Value = 37039892
$ws.Range("K136").Value = 9290.117400000001
$ws.Range("L136").Value = 111119676
$ws.Range("M136").Value = -6740.117400000001
$ws.Range("N136").Value = -111124776
